$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = [double]"11.05753622199721"
$ws.Range("C2").Value = [double]"1.110223024625157e-16"
$ws.Range("D2").Value = [double]"0.005421193709372994"
$ws.Range("E2").Value = [double]"0.7303306206638861"
$ws.Range("F2").Value = [double]"0.5333828154792971"

$ws.Range("B3").Value = [double]"11.30900252532362"
$ws.Range("C3").Value = [double]"1.110223024625157e-16"
$ws.Range("D3").Value = [double]"0.005422450886597523"
$ws.Range("E3").Value = [double]"0.7304999846585886"
$ws.Range("F3").Value = [double]"0.5336302275861982"

$ws.Range("B4").Value = [double]"10.94425829463008"
$ws.Range("C4").Value = [double]"1.110223024625157e-16"
$ws.Range("D4").Value = [double]"0.00539328069052422"
$ws.Range("E4").Value = [double]"0.7265702436190331"
$ws.Range("F4").Value = [double]"0.527904318912621"

$ws.Range("B5").Value = [double]"11.36532305674047"
$ws.Range("C5").Value = [double]"1.110223024625157e-16"
$ws.Range("D5").Value = [double]"0.005564303281522381"
$ws.Range("E5").Value = [double]"0.7496100097161716"
$ws.Range("F5").Value = [double]"0.561915166666679"

$ws.Range("B6").Value = [double]"11.19464679770002"
$ws.Range("C6").Value = [double]"1.110223024625157e-16"
$ws.Range("D6").Value = [double]"0.005394290914333699"
$ws.Range("E6").Value = [double]"0.7267063386234435"
$ws.Range("F6").Value = [double]"0.528102102595491"
